# Delete the SAMSUNG / Charger / Micro USB product row (row 7).
# This consolidates the old delete_product/delete_beverage logic into a
# single delete_item operation: remove the target row entirely and let
# Excel shift the remaining rows up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(7).Delete()
